$d = $word.ActiveDocument

# Locate the "BATMAM" list paragraph.
$batmanPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "BATMAM") {
        $batmanPara = $p
    }
}

# Insert a new paragraph right after it; it inherits the same list
# style/numbering (PargrafodaLista, numId 1) since it follows a list item.
$newRange = $batmanPara.Range.InsertParagraphAfter()

$insertedPara = $batmanPara.Next()
$insertedPara.Range.Text = "Ainda estou aqui "
